$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column (C) for rows 2-11
# from serial date 45170 to 45174, keeping existing formatting.
$ws.Range("C2:C11").Value = 45174
